$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.796.23"
Set-TextValue $ws.Range("E2") "  +0.64%  "
Set-TextValue $ws.Range("D3") "2.477.67"
Set-TextValue $ws.Range("E3") "  +0.30%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  +0.17%  "
Set-TextValue $ws.Range("D5") "319.22"
Set-TextValue $ws.Range("E5") "  +1.45%  "
Set-TextValue $ws.Range("D6") "93.18"
Set-TextValue $ws.Range("E6") "  +1.30%  "
Set-TextValue $ws.Range("E7") "  +0.32%  "
Set-TextValue $ws.Range("E8") "  +0.09%  "
Set-TextValue $ws.Range("D9") "0.518"
Set-TextValue $ws.Range("E9") "  +0.31%  "
Set-TextValue $ws.Range("D10") "0.0876"
Set-TextValue $ws.Range("E10") "  +10.50%  "
Set-TextValue $ws.Range("D11") "33.35"
Set-TextValue $ws.Range("E11") "  +2.42%  "
Set-TextValue $ws.Range("E12") "  +0.69%  "
Set-TextValue $ws.Range("D13") "2.856.88"
Set-TextValue $ws.Range("E13") "  +0.23%  "
Set-TextValue $ws.Range("E14") "  +0.85%  "
Set-TextValue $ws.Range("D15") "15.70"
Set-TextValue $ws.Range("E15") "  -1.89%  "
Set-TextValue $ws.Range("D16") "2.471.60"
Set-TextValue $ws.Range("E16") "  -0.22%  "
Set-TextValue $ws.Range("E17") "  +2.46%  "
Set-TextValue $ws.Range("D18") "41.736.12"
Set-TextValue $ws.Range("E18") "  +0.41%  "
Set-TextValue $ws.Range("D19") "6.48"
Set-TextValue $ws.Range("E19") "  -0.61%  "
Set-TextValue $ws.Range("D20") "0.0₃0951"
Set-TextValue $ws.Range("E20") "  +0.70%  "
Set-TextValue $ws.Range("D21") "71.19"
Set-TextValue $ws.Range("E21") "  +0.10%  "
Set-TextValue $ws.Range("D22") "11.34"
Set-TextValue $ws.Range("E22") "  +1.88%  "
Set-TextValue $ws.Range("D23") "240.90"
Set-TextValue $ws.Range("E23") "  +1.50%  "
Set-TextValue $ws.Range("E24") "  +1.21%  "
Set-TextValue $ws.Range("E25") "  +2.12%  "
Set-TextValue $ws.Range("E26") "  +0.04%  "
Set-TextValue $ws.Range("D27") "24.84"
Set-TextValue $ws.Range("E27") "  +0.28%  "
Set-TextValue $ws.Range("D28") "2.26"
Set-TextValue $ws.Range("E28") "  +0.86%  "
Set-TextValue $ws.Range("D29") "9.79"
Set-TextValue $ws.Range("E29") "  +0.72%  "
Set-TextValue $ws.Range("D30") "36.67"
Set-TextValue $ws.Range("E30") "  +3.24%  "
Set-TextValue $ws.Range("D31") "158.04"
Set-TextValue $ws.Range("E31") "  +1.45%  "
Set-TextValue $ws.Range("D32") "5.51"
Set-TextValue $ws.Range("E32") "  +0.98%  "
Set-TextValue $ws.Range("E33") "  +0.01%  "
Set-TextValue $ws.Range("D34") "0.0766"
Set-TextValue $ws.Range("E34") "  +0.83%  "
Set-TextValue $ws.Range("E35") "  +0.04%  "
Set-TextValue $ws.Range("D36") "17.54"
Set-TextValue $ws.Range("E36") "  +1.34%  "
Set-TextValue $ws.Range("E37") "  +4.47%  "
Set-TextValue $ws.Range("D38") "2.93"
Set-TextValue $ws.Range("E38") "  +1.12%  "
Set-TextValue $ws.Range("E39") "  +1.76%  "
Set-TextValue $ws.Range("E40") "  +0.56%  "
Set-TextValue $ws.Range("D41") "2.54"
Set-TextValue $ws.Range("E41") "  +7.47%  "
Set-TextValue $ws.Range("E42") "  +0.39%  "
Set-TextValue $ws.Range("D43") "2.003.27"
Set-TextValue $ws.Range("E43") "  +2.93%  "
Set-TextValue $ws.Range("D46") "3.00"
Set-TextValue $ws.Range("E46") "  +2.93%  "
Set-TextValue $ws.Range("D47") "9.53"
Set-TextValue $ws.Range("E47") "  +4.67%  "
Set-TextValue $ws.Range("D48") "2.712.76"
Set-TextValue $ws.Range("E48") "  +0.14%  "
Set-TextValue $ws.Range("D49") "98.51"
Set-TextValue $ws.Range("E49") "  +1.15%  "
Set-TextValue $ws.Range("D50") "74.89"
Set-TextValue $ws.Range("E50") "  +4.52%  "
Set-TextValue $ws.Range("D51") "67.37"
Set-TextValue $ws.Range("E51") "  +0.06%  "

# Row 44/45 swap: VeChain/EnergySwap order flips, with new D/E values
Set-TextValue $ws.Range("B44") "EnergySwap"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "19.10"
Set-TextValue $ws.Range("E44") "  +0.61%  "

Set-TextValue $ws.Range("B45") "VeChain"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D45") "0.0285"
Set-TextValue $ws.Range("E45") "  +0.70%  "

Write-Host "Applied cryptos list update."
